# This script applies the "items" workbook edits described in the commit:
#   - renamed spell color to spell type  (header text + key/legend text)
#   - added key to spell slots & spell beans (new "slot 1..5" columns, AC:AG)
#   - drew new icons for equipment items (new icon-sheet coordinates, col D)
#   - preparing for weapon spell slots... (new "weapon chop ..." columns, Z:AB)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): rename the level-item bounding box header and add
#    the new weapon-chop / spell-slot headers in columns Z:AG.
# ---------------------------------------------------------------------------
$ws.Range("V1").Value = "level item bounding box (0, 0, w, h)"

$ws.Range("Z1").Value  = "weapon chop cooldown (ms)"
$ws.Range("AA1").Value = "weapon chop rectangle (0, 0, w, h)"
$ws.Range("AB1").Value = "weapon chop damage"
$ws.Range("AC1").Value = "slot 1 (the first entry is the type (elemental, twilight, necromancy, divine, illusion) and after this, the 0 and 1 show if there is a modifier for this type (strength, duration, range, speed, damage, count, reflect)"
$ws.Range("AD1").Value = "slot 2"
$ws.Range("AE1").Value = "slot 3"
$ws.Range("AF1").Value = "slot 4"
$ws.Range("AG1").Value = "slot 5"

# ---------------------------------------------------------------------------
# 2. New equipment icon-sheet coordinates (column D), rows 2-6 - new icons
#    were drawn on the sprite sheet so their source rectangles moved.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "0, 150"
$ws.Range("D3").Value = "50, 150"
$ws.Range("D4").Value = "200, 150"
$ws.Range("D5").Value = "150, 150"
$ws.Range("D6").Value = "100, 150"

# ---------------------------------------------------------------------------
# 3. Weapon spell-slot prep data for the two weapons (ice staff / rusty
#    sword): chop cooldown, chop rectangle, chop damage, and the two
#    currently-filled spell slots.
# ---------------------------------------------------------------------------
$ws.Range("Z4").Value  = 500
$ws.Range("AA4").Value = "40, 80"
$ws.Range("AB4").Value = 10
$ws.Range("AC4").Value = "1, 0, 0, 0, 1, 0, 1, 0"
$ws.Range("AD4").Value = "1, 0, 0, 0, 0, 1, 0, 1"

$ws.Range("Z5").Value  = 350
$ws.Range("AA5").Value = "40, 50"
$ws.Range("AB5").Value = 5

# ---------------------------------------------------------------------------
# 4. Level-item bounding boxes (column V, rows 7-17) - simplified from the
#    redundant "0, 0, w, h" form down to just "w, h" (matching the new
#    shorter header).
# ---------------------------------------------------------------------------
$ws.Range("V7").Value  = "22, 22"
$ws.Range("V8").Value  = "30, 5"
$ws.Range("V9").Value  = "18, 30"
$ws.Range("V10").Value = "50, 70"
$ws.Range("V11").Value = "30, 30"
$ws.Range("V12").Value = "30, 50"
$ws.Range("V13").Value = "50, 30"
$ws.Range("V14").Value = "30, 40"
$ws.Range("V15").Value = "5, 5"
$ws.Range("V16").Value = "25, 25"
$ws.Range("V17").Value = "5, 5"
